# Applies the cryptos list price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, preserving the exact string
# (e.g. "5.930", "8.200") even when it would otherwise parse as a number,
# and restore the cells original style afterward so no formatting changes
# are introduced.
function Set-TextValue {
    param($Worksheet, $CellRef, $Text)
    $cell = $Worksheet.Range($CellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value2 = $Text
    $cell.Style = $origStyle
}

$ws.Range('D2').Value2 = '27.538.07'
$ws.Range('E2').Value2 = '  +5.46%  '
$ws.Range('D3').Value2 = '1.726.28'
$ws.Range('E3').Value2 = '  +4.52%  '
Set-TextValue $ws 'D5' '225.67'
$ws.Range('E5').Value2 = '  +3.30%  '
Set-TextValue $ws 'D6' '0.5366'
$ws.Range('E6').Value2 = '  +2.89%  '
$ws.Range('E7').Value2 = '  +0.01%  '
Set-TextValue $ws 'D8' '0.2674'
$ws.Range('E8').Value2 = '  +1.08%  '
Set-TextValue $ws 'D9' '0.06605'
$ws.Range('E9').Value2 = '  +4.23%  '
Set-TextValue $ws 'D10' '21.79'
$ws.Range('E10').Value2 = '  +6.78%  '
Set-TextValue $ws 'D11' '0.07705'
$ws.Range('E11').Value2 = '  +0.12%  '
$ws.Range('E12').Value2 = '  -0.27%  '
$ws.Range('D13').Value2 = '1.722.40'
$ws.Range('E13').Value2 = '  +3.93%  '
$ws.Range('D14').Value2 = '1.963.39'
$ws.Range('E14').Value2 = '  +4.46%  '
Set-TextValue $ws 'D15' '0.5851'
$ws.Range('E15').Value2 = '  +4.59%  '
$ws.Range('E16').Value2 = '  +1.84%  '
Set-TextValue $ws 'D17' '68.06'
$ws.Range('E17').Value2 = '  +4.05%  '
$ws.Range('D18').Value2 = '27.562.70'
$ws.Range('E18').Value2 = '  +5.55%  '
Set-TextValue $ws 'D19' '221.07'
$ws.Range('E19').Value2 = '  +15.58%  '
$ws.Range('E20').Value2 = '  +0.09%  '
$ws.Range('E21').Value2 = '  +2.25%  '
Set-TextValue $ws 'D22' '10.66'
$ws.Range('E22').Value2 = '  +1.73%  '
Set-TextValue $ws 'D23' '6.098'
$ws.Range('E23').Value2 = '  +2.76%  '
Set-TextValue $ws 'D24' '1.004'
$ws.Range('E24').Value2 = '  +0.03%  '
Set-TextValue $ws 'D25' '148.45'
$ws.Range('E25').Value2 = '  +2.10%  '
Set-TextValue $ws 'D26' '1.713'
$ws.Range('E26').Value2 = '  +13.38%  '
Set-TextValue $ws 'D27' '0.1234'
$ws.Range('E27').Value2 = '  +3.66%  '
Set-TextValue $ws 'D28' '7.411'
$ws.Range('E28').Value2 = '  +2.58%  '
Set-TextValue $ws 'D29' '16.69'
Set-TextValue $ws 'D30' '0.05577'
$ws.Range('E30').Value2 = '  +2.02%  '
$ws.Range('E31').Value2 = '  +2.46%  '
Set-TextValue $ws 'D32' '3.557'
$ws.Range('E32').Value2 = '  +3.16%  '
Set-TextValue $ws 'D33' '3.460'
$ws.Range('E33').Value2 = '  +2.85%  '
Set-TextValue $ws 'D34' '1.659'
$ws.Range('E34').Value2 = '  +6.49%  '
Set-TextValue $ws 'D35' '0.9636'
$ws.Range('E35').Value2 = '  +1.31%  '
Set-TextValue $ws 'D36' '2.825'
$ws.Range('E36').Value2 = '  +1.42%  '
Set-TextValue $ws 'D37' '2.431'
$ws.Range('E37').Value2 = '  +1.25%  '
Set-TextValue $ws 'D38' '0.5953'
$ws.Range('E38').Value2 = '  +5.59%  '
Set-TextValue $ws 'D39' '0.01649'
$ws.Range('E39').Value2 = '  +4.42%  '
Set-TextValue $ws 'D40' '5.930'
$ws.Range('E40').Value2 = '  +1.20%  '
Set-TextValue $ws 'D41' '0.8568'
$ws.Range('E41').Value2 = '  +2.94%  '
$ws.Range('D42').Value2 = '1.056.76'
$ws.Range('E42').Value2 = '  +2.67%  '
$ws.Range('E43').Value2 = '  +0.05%  '
Set-TextValue $ws 'D44' '101.48'
$ws.Range('E44').Value2 = '  +0.24%  '
$ws.Range('D45').Value2 = '1.870.27'
$ws.Range('E45').Value2 = '  +4.40%  '
$ws.Range('E46').Value2 = '  +6.31%  '
Set-TextValue $ws 'D47' '59.09'
$ws.Range('E47').Value2 = '  +2.47%  '
Set-TextValue $ws 'D48' '8.200'
$ws.Range('E48').Value2 = '  +2.55%  '
Set-TextValue $ws 'D49' '0.4436'
$ws.Range('E49').Value2 = '  +2.24%  '
$ws.Range('E50').Value2 = '  +0.24%  '
Set-TextValue $ws 'D51' '0.05265'
$ws.Range('E51').Value2 = '  +1.67%  '
